$d = $word.ActiveDocument

$replacements = @(
    @("391×2=", "775×8="),
    @("462×2=", "960×3="),
    @("259×6=", "701×8="),
    @("674×4=", "705×8="),
    @("631×3=", "201×5="),
    @("735×2=", "607×3="),
    @("904×5=", "121×4="),
    @("295×3=", "573×3="),
    @("762×5=", "266×7="),
    @("122×3=", "655×6="),
    @("613×3=", "636×5="),
    @("733×5=", "366×6="),
    @("113×3=", "516×4="),
    @("710×5=", "799×6="),
    @("188×4=", "780×4="),
    @("586×4=", "346×3="),
    @("138×2=", "160×5="),
    @("780×6=", "658×9="),
    @("945×8=", "426×3="),
    @("204×9=", "416×5="),
    @("738×4=", "933×3="),
    @("443×4=", "640×5="),
    @("565×3=", "184×3="),
    @("806×2=", "202×9="),
    @("209×3=", "929×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
